$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials-x-IMU3-SA-Syn")

# Cell C6 holds the Comment "180R" for designators R1, R2 (Resistor 0402).
# Replace it with "470R", preserving the existing text formatting
# (leading apostrophe keeps Excel's quote-prefix / text style, matching
# the original cell's formatting).
$ws.Range("C6").Value = "'470R"
